# Correction in SA algorithm run_28 fitness log (rows reflect corrected
# "best fitness so far" values after fixing the simulated annealing run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17 (Generation 0-15): corrected best-fitness value -> 7734
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 3).Value = 7734
}

# Rows 18-124 (Generation 16-122): corrected best-fitness value -> 7310
for ($r = 18; $r -le 124; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}

# Rows 125-164 (Generation 123-162): corrected best-fitness value -> 7293
for ($r = 125; $r -le 164; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}

Write-Host "Applied SA algorithm correction to run_28 fitness log (rows 2-164)."
